$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) / Volume(1h) (E) columns with the latest scrape.
# Some D values are plain decimals that the COM layer would auto-coerce
# to numbers (losing trailing zeros / using binary-float repr like
# "1.7399999999999999E-5" instead of "0.0000174"); force those cells to
# Text first, then drop the style back to Normal so no stray number
# format lingers on the cell.
$ws.Range("D2").Value = "69.403.76"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.484.44"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.484.21"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "2.944.45"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "69.337.49"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "2.483.20"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.99%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").Value = "2.621.45"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "0.0₃0864"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "435.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.313"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.70%  "
$ws.Range("E49").Value = "  -1.39%  "

# Rows 50/51 swap: POPCAT moves up to row 50 (now up +22.54%), and
# Mantle drops to row 51 with a refreshed price/volume.
$ws.Range("B50").Value = "POPCAT"
$ws.Range("C50").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +22.54%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.571"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "
